$wb = $excel.ActiveWorkbook

# "Active" sheet (sheet1): the "100%" size button task is complete, remove it from Active
$active = $wb.Worksheets.Item("Active")
$active.Rows("7:7").Delete()

# "Inactive" sheet (sheet2): add the completed task as a new row at the top of the data
$inactive = $wb.Worksheets.Item("Inactive")
$inactive.Rows("2:2").Insert()

# the insert copies the bold header formatting down; reset the new row back to normal
$inactive.Range("A2:F2").Font.Bold = $false

$inactive.Cells.Item(2, 1).Value = 46
$inactive.Cells.Item(2, 2).Value = "add ""100%"" size button"
$inactive.Cells.Item(2, 3).Value = "Done"
$inactive.Cells.Item(2, 4).Value = "Task"
$inactive.Cells.Item(2, 5).Value = "'8/16/2018"
$inactive.Cells.Item(2, 6).Value = "'8/21/2018"
